# 70. Gradient Descent in Action
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the old "Price"/"Sqrt" header row (A1:B1) down to A9:B9 ---
# (this is the old "Guess for B" / "MSE" table's header, which shifts down
#  while the top header becomes "Lot Size" / "House price")
$ws.Range("A9").Value = "Guess for B"
$ws.Range("B9").Value = "MSE"

# --- Replace the top header with the new labels ---
$ws.Range("A1").Value = "Lot Size"
$ws.Range("B1").Value = "House price"

# --- Bold the two header rows ---
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A9:B9").Font.Bold = $true

# --- Gradient descent scratch area (columns G, I, J, K, L) ---
$ws.Range("G2").Value = "Learning Rate"
$ws.Range("G2").Font.Bold = $true
$ws.Range("G3").Value = 0.3

$ws.Range("I2").Value = "Initial Guess:"
$ws.Range("I2").Font.Bold = $true

$ws.Range("J1").Value = "Value of B"
$ws.Range("J1").Font.Bold = $true
$ws.Range("K1").Value = "Slope of MSE"
$ws.Range("K1").Font.Bold = $true
$ws.Range("L1").Value = "How Much To Adjust B"
$ws.Range("L1").Font.Bold = $true

# Initial guess value
$ws.Range("J2").Value = 0

# Slope-of-MSE and adjustment formulas, row 2
$ws.Range("K2").Formula = "=2 * ((J2-`$B`$2) + (J2-`$B`$3) + (J2-`$B`$4) + (J2-`$B`$5) + (J2-`$B`$6) + (J2-`$B`$7)) / 6"
$ws.Range("L2").Formula = "=K2*`$G`$3"

# Next guess (row 3) references the previous row
$ws.Range("J3").Formula = "=J2-L2"
$ws.Range("K3").Formula = "=2 * ((J3-`$B`$2) + (J3-`$B`$3) + (J3-`$B`$4) + (J3-`$B`$5) + (J3-`$B`$6) + (J3-`$B`$7)) / 6"
$ws.Range("L3").Formula = "=K3*`$G`$3"

# Fill J/K/L down through row 41, letting Excel auto-adjust relative refs
$ws.Range("J4:J41").Formula = "=J3-L3"
$ws.Range("K4:K41").Formula = "=2 * ((J4-`$B`$2) + (J4-`$B`$3) + (J4-`$B`$4) + (J4-`$B`$5) + (J4-`$B`$6) + (J4-`$B`$7)) / 6"
$ws.Range("L4:L41").Formula = "=K4*`$G`$3"

# K column wraps text
$ws.Range("K2:K41").WrapText = $true

# Highlight the row where B has essentially converged (row 24)
$ws.Range("J24:L24").Interior.Color = 65535
$ws.Range("K24").WrapText = $true

# Column K width
$ws.Range("K1").EntireColumn.AutoFit()

$ws.Range("M3").Select()
